$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 13:50"

# Row 18/19: Portugal overtakes Canada in the ranking, so the two rows swap
# places; Portugal's figures are also refreshed with newer totals while
# Canada's row keeps the exact same figures it had before.
$ws.Cells.Item(18,1).Value = "Portugal"
$ws.Cells.Item(18,2).Value = 5962
$ws.Cells.Item(18,3).Value = 792
$ws.Cells.Item(18,4).Value = 43
$ws.Cells.Item(18,5).Value = 5800
$ws.Cells.Item(18,6).Value = 89
$ws.Cells.Item(18,7).Value = 19
$ws.Cells.Item(18,8).Value = 119

$ws.Cells.Item(19,1).Value = "Canada"
$ws.Cells.Item(19,2).Value = 5655
$ws.Cells.Item(19,3).Value = 0
$ws.Cells.Item(19,4).Value = 508
$ws.Cells.Item(19,5).Value = 5087
$ws.Cells.Item(19,6).Value = 120
$ws.Cells.Item(19,7).Value = 0
$ws.Cells.Item(19,8).Value = 60

# Row 20 (Noruega): refreshed totals
$ws.Cells.Item(20,2).Value = 4213
$ws.Cells.Item(20,3).Value = 198
$ws.Cells.Item(20,5).Value = 4183

# Row 22 (Brasil): refreshed totals
$ws.Cells.Item(22,5).Value = 3781
$ws.Cells.Item(22,7).Value = 3
$ws.Cells.Item(22,8).Value = 117

# Row 76 (Principado de Andorra): refreshed totals
$ws.Cells.Item(76,5).Value = 303
$ws.Cells.Item(76,7).Value = 1
$ws.Cells.Item(76,8).Value = 4

# Rows 143-147: Bermudas overtakes Etiopia/Niger/Mali/Islas Virgenes de los
# Estados Unidos in the ranking, so it moves up to row 143 (with refreshed
# totals) and the other four rows shift down by one, keeping their figures.
$ws.Cells.Item(143,1).Value = "Bermudas"
$ws.Cells.Item(143,2).Value = 22
$ws.Cells.Item(143,3).Value = 5
$ws.Cells.Item(143,4).Value = 2
$ws.Cells.Item(143,5).Value = 20
$ws.Cells.Item(143,6).Value = 0
$ws.Cells.Item(143,7).Value = 0
$ws.Cells.Item(143,8).Value = 0

$ws.Cells.Item(144,1).Value = "Etiopia"
$ws.Cells.Item(144,2).Value = 19
$ws.Cells.Item(144,3).Value = 3
$ws.Cells.Item(144,4).Value = 1
$ws.Cells.Item(144,5).Value = 18
$ws.Cells.Item(144,6).Value = 0
$ws.Cells.Item(144,7).Value = 0
$ws.Cells.Item(144,8).Value = 0

$ws.Cells.Item(145,1).Value = "Niger"
$ws.Cells.Item(145,2).Value = 18
$ws.Cells.Item(145,3).Value = 8
$ws.Cells.Item(145,4).Value = 0
$ws.Cells.Item(145,5).Value = 17
$ws.Cells.Item(145,6).Value = 0
$ws.Cells.Item(145,7).Value = 0
$ws.Cells.Item(145,8).Value = 1

$ws.Cells.Item(146,1).Value = "Mali"
$ws.Cells.Item(146,2).Value = 18
$ws.Cells.Item(146,3).Value = 0
$ws.Cells.Item(146,4).Value = 0
$ws.Cells.Item(146,5).Value = 17
$ws.Cells.Item(146,6).Value = 0
$ws.Cells.Item(146,7).Value = 0
$ws.Cells.Item(146,8).Value = 1

$ws.Cells.Item(147,1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(147,2).Value = 17
$ws.Cells.Item(147,3).Value = 0
$ws.Cells.Item(147,4).Value = 0
$ws.Cells.Item(147,5).Value = 17
$ws.Cells.Item(147,6).Value = 0
$ws.Cells.Item(147,7).Value = 0
$ws.Cells.Item(147,8).Value = 0
